# This sheet stores a rolling buffer of accelerometer/gyroscope samples
# (columns: timestamp, label, ax, ay, az, gx, gy, gz).
# On 2023-05-09 a fresh batch of 18 samples was captured and prepended to
# the buffer; the buffer is capped at 30 data rows, so the oldest samples
# fall off the end. The final state (30 data rows + header) is written
# directly below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 30,8
$data[0,0] = 0
$data[0,1] = "falling"
$data[0,2] = -3.58582592010498
$data[0,3] = 8.237812042236328
$data[0,4] = 0.2284512519836425
$data[0,5] = -0.0005742134153842587
$data[0,6] = 0.0210748501121997
$data[0,7] = 0.04580267548561091
$data[1,0] = 100
$data[1,1] = "falling"
$data[1,2] = -3.345468521118164
$data[1,3] = 8.276437759399414
$data[1,4] = -0.0355764925479888
$data[1,5] = -0.0293826170265674
$data[1,6] = 0.0331699833273887
$data[1,7] = 0.01545489076524965
$data[2,0] = 200
$data[2,1] = "falling"
$data[2,2] = -3.182104587554932
$data[2,3] = 8.366311073303223
$data[2,4] = 0.2440185844898224
$data[2,5] = -0.02163684628903856
$data[2,6] = 0.01164309203624722
$data[2,7] = -0.01903456177562478
$data[3,0] = 300
$data[3,1] = "falling"
$data[3,2] = -3.532341480255127
$data[3,3] = 8.214254379272461
$data[3,4] = 0.022090196609497
$data[3,5] = -0.00761748962104312
$data[3,6] = -0.02296242564916602
$data[3,7] = -0.00737925238907336
$data[4,0] = 400
$data[4,1] = "falling"
$data[4,2] = -3.395359516143799
$data[4,3] = 8.257164001464844
$data[4,4] = 0.0318913161754608
$data[4,5] = -0.03987117595970628
$data[4,6] = -0.01981036089360703
$data[4,7] = 0.01966986127197736
$data[5,0] = 500
$data[5,1] = "falling"
$data[5,2] = -3.210868835449219
$data[5,3] = 8.186161994934082
$data[5,4] = -0.0299544632434844
$data[5,5] = -0.04132503550499657
$data[5,6] = -0.01673159938305617
$data[5,7] = 0.0125227374956011
$data[6,0] = 600
$data[6,1] = "falling"
$data[6,2] = -3.197915077209473
$data[6,3] = 8.212710380554199
$data[6,4] = 0.0095521211624145
$data[6,5] = -0.01340849157422786
$data[6,6] = -0.03075706511735912
$data[6,7] = -0.009834930114448132
$data[7,0] = 700
$data[7,1] = "falling"
$data[7,2] = -3.255885124206543
$data[7,3] = 8.21357250213623
$data[7,4] = -0.070087194442749
$data[7,5] = -0.03994447708129877
$data[7,6] = -0.01221119597554203
$data[7,7] = -0.04150218397378917
$data[8,0] = 800
$data[8,1] = "falling"
$data[8,2] = -3.168186187744141
$data[8,3] = 8.229218482971191
$data[8,4] = -0.0939638018608093
$data[8,5] = -0.03637702405452719
$data[8,6] = 0.01701259657740584
$data[8,7] = -0.03676186949014657
$data[9,0] = 900
$data[9,1] = "falling"
$data[9,2] = -3.290424346923828
$data[9,3] = 8.159111976623535
$data[9,4] = -0.1956053972244262
$data[9,5] = -0.03740938737988467
$data[9,6] = 0.01381166309118263
$data[9,7] = -0.03841120541095729
$data[10,0] = 1000
$data[10,1] = "falling"
$data[10,2] = -3.586246490478516
$data[10,3] = 8.05996036529541
$data[10,4] = -0.0569053888320922
$data[10,5] = -0.04191757388412946
$data[10,6] = 0.01410487815737716
$data[10,7] = -0.03918089691549535
$data[11,0] = 1100
$data[11,1] = "falling"
$data[11,2] = -3.300580024719238
$data[11,3] = 8.124805450439453
$data[11,4] = -0.1510338187217712
$data[11,5] = 0.01087340153753754
$data[11,6] = 0.01499674115329971
$data[11,7] = -0.01282817013561728
$data[12,0] = 1200
$data[12,1] = "falling"
$data[12,2] = -2.94456958770752
$data[12,3] = 8.213338851928711
$data[12,4] = -0.0691232085227966
$data[12,5] = 0.02458121769130225
$data[12,6] = 0.01081842321902504
$data[12,7] = -0.04899139240384098
$data[13,0] = 1300
$data[13,1] = "falling"
$data[13,2] = -3.094478130340576
$data[13,3] = 8.185011863708496
$data[13,4] = -0.1414701342582702
$data[13,5] = 0.01802052438259117
$data[13,6] = 0.00392786357551804
$data[13,7] = -0.03619987547397605
$data[14,0] = 1400
$data[14,1] = "falling"
$data[14,2] = -3.334782123565674
$data[14,3] = 8.0909423828125
$data[14,4] = -0.0672928094863891
$data[14,5] = 0.01199739351868624
$data[14,6] = -0.0002138027921320222
$data[14,7] = 0.03830125063657756
$data[15,0] = 1500
$data[15,1] = "falling"
$data[15,2] = -3.368669509887696
$data[15,3] = 8.071453094482422
$data[15,4] = 0.1171565353870391
$data[15,5] = 0.01877188928425304
$data[15,6] = -0.01797165483236304
$data[15,7] = 0.03120299618691198
$data[16,0] = 1600
$data[16,1] = "falling"
$data[16,2] = -3.412579536437988
$data[16,3] = 8.035589218139648
$data[16,4] = 0.0592367351055145
$data[16,5] = -0.0108428578823804
$data[16,6] = -0.0058032199740409
$data[16,7] = -0.0042760567739605
$data[17,0] = 1700
$data[17,1] = "falling"
$data[17,2] = -3.019937038421631
$data[17,3] = 8.054733276367188
$data[17,4] = -0.1420263051986694
$data[17,5] = 0.01026864476501944
$data[17,6] = 0.01618792921304707
$data[17,7] = 0.006133087240159481
$data[18,0] = 1800
$data[18,1] = "falling"
$data[18,2] = -3.012916564941406
$data[18,3] = 8.089370727539062
$data[18,4] = -0.1633265316486358
$data[18,5] = 0.03629761248826986
$data[18,6] = 0.01907121278345579
$data[18,7] = 0.05546045627444995
$data[19,0] = 1900
$data[19,1] = "falling"
$data[19,2] = -3.395848751068115
$data[19,3] = 8.023316383361816
$data[19,4] = 0.0382503271102905
$data[19,5] = 0.05165476366877556
$data[19,6] = -0.0003787364251911958
$data[19,7] = 0.03377473920583711
$data[20,0] = 2000
$data[20,1] = "falling"
$data[20,2] = -3.384797096252441
$data[20,3] = 7.934267520904541
$data[20,4] = 0.07479587197303771
$data[20,5] = 0.02702467799186697
$data[20,6] = -0.02729956846684218
$data[20,7] = -0.006963863894343374
$data[21,0] = 2100
$data[21,1] = "falling"
$data[21,2] = -3.632324695587158
$data[21,3] = 7.965863227844238
$data[21,4] = 0.0220168232917785
$data[21,5] = 0.02345722466707222
$data[21,6] = -0.01078177168965329
$data[21,7] = -0.01979203335940831
$data[22,0] = 2200
$data[22,1] = "falling"
$data[22,2] = -3.280028343200684
$data[22,3] = 7.839587211608887
$data[22,4] = 0.0374422371387481
$data[22,5] = 0.02312735825777049
$data[22,6] = -0.06875288158655167
$data[22,7] = -0.001032362207770397
$data[23,0] = 2300
$data[23,1] = "falling"
$data[23,2] = -3.331558704376221
$data[23,3] = 7.833842277526855
$data[23,4] = -0.1311583817005157
$data[23,5] = -0.1697350136935719
$data[23,6] = -0.2189157873392111
$data[23,7] = -0.04578435219824339
$data[24,0] = 2400
$data[24,1] = "falling"
$data[24,2] = -3.715910911560059
$data[24,3] = 7.601772308349609
$data[24,4] = -0.8390151262283325
$data[24,5] = -0.2920058012008669
$data[24,6] = -0.2727635514736176
$data[24,7] = 0.126033713221551
$data[25,0] = 2500
$data[25,1] = "falling"
$data[25,2] = -3.635088920593262
$data[25,3] = 7.303267478942871
$data[25,4] = -1.838643550872803
$data[25,5] = -0.2655186891555784
$data[25,6] = -0.2190562760829923
$data[25,7] = 0.6186354464292545
$data[26,0] = 2600
$data[26,1] = "falling"
$data[26,2] = -3.690509796142578
$data[26,3] = 6.705544471740723
$data[26,4] = -2.344899177551269
$data[26,5] = -0.2926044583320622
$data[26,6] = -0.3912286460399645
$data[26,7] = 0.472412636876103
$data[27,0] = 2700
$data[27,1] = "falling"
$data[27,2] = -4.218070983886719
$data[27,3] = 5.676623344421387
$data[27,4] = -2.678865909576416
$data[27,5] = -0.02654820919036799
$data[27,6] = 0.1870408368110673
$data[27,7] = -0.7209736722707767
$data[28,0] = 2800
$data[28,1] = "falling"
$data[28,2] = -5.941818714141846
$data[28,3] = 3.218802928924561
$data[28,4] = -3.337108850479126
$data[28,5] = 1.565036740303041
$data[28,6] = -0.7936239337921156
$data[28,7] = -2.525304698944093
$data[29,0] = 2900
$data[29,1] = "falling"
$data[29,2] = -7.670186996459961
$data[29,3] = 2.646533727645874
$data[29,4] = -3.358609437942505
$data[29,5] = 1.211272468566887
$data[29,6] = -0.4785640525817794
$data[29,7] = -2.630715656280514

$lastRow = 1 + 30
$ws.Range("A2:H$lastRow").Value = $data

Write-Output "Updated data rows 2..$lastRow"
